# Apply cell updates per diff (cryptos.xlsx refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns that hold numeric-looking price text stay as Text so
# Excel does not silently convert "1.00" -> 1 or "0.125" -> 0.125 (number).
$ws.Range("D2").Value = "41.753.76"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.225.41"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.35"
$ws.Range("E5").Value = "  +6.63%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.82"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  +7.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.04"
$ws.Range("E10").Value = "  +12.20%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.47"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.23"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "2.556.04"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.97"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.863"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "2.230.11"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "41.729.69"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.91"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.10"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +6.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("E25").Value = "  +9.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.53"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.74"
$ws.Range("E28").Value = "  +7.53%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.28"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.80"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.60"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.72"
$ws.Range("E37").Value = "  +16.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.01"
$ws.Range("E38").Value = "  +9.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0303"
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.60"
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.22"
$ws.Range("E43").Value = "  +18.94%  "
$ws.Range("E44").Value = "  +6.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.87"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.77"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.66"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +6.35%  "
$ws.Range("E51").Value = "  +0.33%  "
